# Perbaikan Semua (profile, stok, detail, detail penjualan)
#
# Update the sample "supplier" row (2nd data row) on the "Data supplier"
# sheet: SUP5 / Supplier E / Jl. Melati -> SUP6 / Supplier F / Jl. Bunga,
# and leave the selection on cell C2 (as last used/saved by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data supplier")

$ws.Range("A2").Value = "SUP6"
$ws.Range("B2").Value = "Supplier F"
$ws.Range("C2").Value = "Jl. Bunga"

$ws.Activate()
$ws.Range("C2").Select()

$wb.Save()
